$d = $word.ActiveDocument

# ==========================================================================
# Edit 1: "Values for last year included (e.g. ...)" -- the text itself is
# unchanged, but the three runs that were split around the "e.g." proofing
# marks (gramStart/gramEnd) collapse into a single run. A find/replace over
# the full (identical) visible text merges those runs into one, dropping
# the now-orphaned <w:proofErr/> markers.
# ==========================================================================
$ok1 = $d.Content.Find.Execute(
    "Values for last year included (e.g. relative growth rate measured 2019-2021; 2021 values are below)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Values for last year included (e.g. relative growth rate measured 2019-2021; 2021 values are below)",
    2)

# ==========================================================================
# Edit 2: PVE methods paragraph -- re-describe how PVE was calculated for
# the herbivore/SLA variables (now refit as general linear mixed models)
# instead of the old get_variance()/insight description, and reuse the
# existing "random effect variance/(...) ... VarCorr() ... lme4" text that
# used to describe the "remaining" variables, tacking on a final clause.
# ==========================================================================

# Step 1: drop the now-redundant "...R package. Remaining variables were
# analyzed with general linear mixed models and PVE was calculated as:
# random effect variance/(random effect variance + residual variance) with
# the " sentence sitting between "insight" and "VarCorr" -- its content is
# now said only once.
$ok2 = $d.Content.Find.Execute(
    " R package. Remaining variables were analyzed with general linear mixed models and PVE was calculated as: random effect variance/(random effect variance + residual variance) with the ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2)

# Step 2: turn the now-orphaned italic "insight" run into "VarCorr()" in
# place (keeps that run's own formatting/identity -- it stays italic).
$ok3 = $d.Content.Find.Execute(
    "insight",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "VarCorr()",
    2)

# Step 3: the original "VarCorr()" mention immediately follows and is now a
# duplicate -- collapse "VarCorr()VarCorr()" down to a single "VarCorr()"
# (the merged range keeps the formatting of the first occurrence above).
$ok4 = $d.Content.Find.Execute(
    "VarCorr()VarCorr()",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "VarCorr()",
    2)

# Step 4: reword the lead-in sentence describing how PVE was calculated.
$ok5 = $d.Content.Find.Execute(
    "as: random effect variance/(random effect variance + residual variance) with the get_variance() function from the ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "by refitting the variables to general linear mixed models (i.e., Gaussian distributions), then calculating: random effect variance/(random effect variance + residual variance) with the ",
    2)

# Step 5: append the closing clause after "R package".
$ok6 = $d.Content.Find.Execute(
    " R package.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " R package, as was performed for the other variables.",
    2)

Write-Output "edit1=$ok1 edit2.1=$ok2 edit2.2=$ok3 edit2.3=$ok4 edit2.4=$ok5 edit2.5=$ok6"
